# añadidos 125 y 250 hasta el año 2003
# Adds a new "Class" column (filled with "motogp") to the Equipos team
# table, expanding Table_1 from A1:G40 to A1:H40, and leaves the
# "Equipos" sheet/selection active (matching the authored edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equipos")

# Make this the active sheet (the edit was made here).
$ws.Activate()

$lo = $ws.ListObjects.Item("Table_1")

# Add the new trailing column to the table and name it "Class".
$newCol = $lo.ListColumns.Add()
$ws.Range("H1").Value = "Class"

# Fill every existing data row (2-39) with the class "motogp". Row 40 is
# a trailing blank table row (already reserved by the table's range
# before this edit) and must stay empty.
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 8).Value = "motogp"
}

# Normalize the format of the whole new column (header row already text,
# this also materializes the trailing blank table row beneath the data).
$ws.Range("H2:H40").NumberFormat = "General"

# Match the selection left behind by the edit.
$ws.Range("H2:H39").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
